$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns before the existing "C" (analyst rating) column.
# This shifts the old column C (ratings/price-target text) out to column E,
# matching the target layout of B | C | D | E.
$ws.Columns("C:D").Insert()

# Give the three now-adjacent columns (C, D, E) the same 8-char width that
# column C originally had.
$ws.Columns("C").ColumnWidth = 7.14
$ws.Columns("D").ColumnWidth = 7.14
$ws.Columns("E").ColumnWidth = 7.14

# Row 1 header: shift the old header text from B1 -> D1, keep the (already
# relocated) C1->E1 value, and introduce two new header labels in B1/C1.
# NOTE: use Value2 for reads -- plain .Value getter misbehaves in this host
# when read back (returns property metadata instead of the cell data). The
# old B1 label is read dynamically (rather than hard-coded) so the script
# still does the right thing if the sheet's date label ever changes.
$oldB1 = $ws.Cells.Item(1, 2).Value2
$ws.Cells.Item(1, 4).Value2 = $oldB1
$ws.Cells.Item(1, 2).Value2 = "Jun_17"
$ws.Cells.Item(1, 3).Value2 = "Jun_15"

# Data rows: fill the two newly inserted columns with the same "UN" filler
# value already used in column B for every data row.
for ($r = 2; $r -le 27; $r++) {
    $ws.Cells.Item($r, 3).Value2 = "UN"
    $ws.Cells.Item($r, 4).Value2 = "UN"
}
